$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value (Invoice name) from "Invoice Split-1" to "Default Invoice-3"
$ws.Range("B2").Value = "Default Invoice-3"

# Set column B width
$ws.Columns.Item(2).ColumnWidth = 18.5703125

# Update the selected cell/range to B2
$ws.Range("B2").Select()
